$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A
$ws.Cells.Item(2, 1).Value = "Acordo Paris"
$ws.Cells.Item(3, 1).Value = "Agenda verde"
$ws.Cells.Item(4, 1).Value = "Ambiental"
$ws.Cells.Item(5, 1).Value = "Aquecimento global"
$ws.Cells.Item(6, 1).Value = "Biodiversidade"
$ws.Cells.Item(7, 1).Value = "Camada ozônio"
$ws.Cells.Item(8, 1).Value = "Carbono"
$ws.Cells.Item(9, 1).Value = "Certificação ambiental"
$ws.Cells.Item(10, 1).Value = "Crédito verde"
$ws.Cells.Item(11, 1).Value = "CO2"
$ws.Cells.Item(12, 1).Value = "Combustível limpo"
$ws.Cells.Item(13, 1).Value = "Desmatamento"
$ws.Cells.Item(14, 1).Value = "Descarbonização"
$ws.Cells.Item(15, 1).Value = "Economia verde"
$ws.Cells.Item(16, 1).Value = "Efeito Estufa"
$ws.Cells.Item(17, 1).Value = "Energia limpa"
$ws.Cells.Item(18, 1).Value = "Energia renovável"
$ws.Cells.Item(19, 1).Value = "Energia verde"
$ws.Cells.Item(20, 1).Value = "Gerenciamento resíduos"
$ws.Cells.Item(21, 1).Value = "GHG"
$ws.Cells.Item(22, 1).Value = "Greenwashing"
$ws.Cells.Item(23, 1).Value = "Hidrogênio verde"
$ws.Cells.Item(24, 1).Value = "Meio ambiente"
$ws.Cells.Item(25, 1).Value = "Metano"
$ws.Cells.Item(26, 1).Value = "Mudanças climáticas"
$ws.Cells.Item(27, 1).Value = "Mudanças uso solo"
$ws.Cells.Item(28, 1).Value = "Net Zero"
$ws.Cells.Item(29, 1).Value = "ODS"
$ws.Cells.Item(30, 1).Value = "Poluição"
$ws.Cells.Item(31, 1).Value = "Poluentes"
$ws.Cells.Item(32, 1).Value = "Qualidade ar"
$ws.Cells.Item(33, 1).Value = "Reciclagem"
$ws.Cells.Item(34, 1).Value = "Resíduos tóxicos"
$ws.Cells.Item(35, 1).Value = "Reflorestamento"
$ws.Cells.Item(36, 1).Value = "Risco climático"
$ws.Cells.Item(37, 1).Value = "Usina eólica"
$ws.Cells.Item(38, 1).Value = "Usina solar"
$ws.Cells.Item(39, 1).Value = "Uso racional"

# Column B
$ws.Cells.Item(2, 2).Value = "Bem-estar"
$ws.Cells.Item(3, 2).Value = "Burnout"
$ws.Cells.Item(4, 2).Value = "Desigualdade social"
$ws.Cells.Item(5, 2).Value = "Direitos humanos"
$ws.Cells.Item(6, 2).Value = "Discriminação racial"
$ws.Cells.Item(7, 2).Value = "Diversidade"
$ws.Cells.Item(8, 2).Value = "Doação"
$ws.Cells.Item(9, 2).Value = "Equidade racial"
$ws.Cells.Item(10, 2).Value = "Equidade salarial"
$ws.Cells.Item(11, 2).Value = "Equidade gênero"
$ws.Cells.Item(12, 2).Value = "Escravidão"
$ws.Cells.Item(13, 2).Value = "Filantropia"
$ws.Cells.Item(14, 2).Value = "Gay"
$ws.Cells.Item(15, 2).Value = "Grupos minorizados"
$ws.Cells.Item(16, 2).Value = "Igualdade de gênero"
$ws.Cells.Item(17, 2).Value = "Igualdade racial"
$ws.Cells.Item(18, 2).Value = "Impacto social"
$ws.Cells.Item(19, 2).Value = "Inclusão digital"
$ws.Cells.Item(20, 2).Value = "Inclusão social"
$ws.Cells.Item(21, 2).Value = "Intolerância"
$ws.Cells.Item(22, 2).Value = "LGBT"
$ws.Cells.Item(23, 2).Value = "LGBTQIA"
$ws.Cells.Item(24, 2).Value = "Liderança feminina"
$ws.Cells.Item(25, 2).Value = "Mulheres liderança"
$ws.Cells.Item(26, 2).Value = "Pessoas negras"
$ws.Cells.Item(27, 2).Value = "Pessoas pretas"
$ws.Cells.Item(28, 2).Value = "Preconceito"
$ws.Cells.Item(29, 2).Value = "Racismo"
$ws.Cells.Item(30, 2).Value = "Relações trabalhistas"
$ws.Cells.Item(31, 2).Value = "Responsabilidade social"
$ws.Cells.Item(32, 2).Value = "Trabalho escravo"
$ws.Cells.Item(33, 2).Value = "Trabalho infantil"
$ws.Cells.Item(34, 2).Value = "Transgêneros"
$ws.Cells.Item(35, 2).Value = "Violência"

# Column C
$ws.Cells.Item(2, 3).Value = "Ação penal"
$ws.Cells.Item(3, 3).Value = "Anticorrupção"
$ws.Cells.Item(4, 3).Value = "Assédio moral"
$ws.Cells.Item(5, 3).Value = "Auditoria"
$ws.Cells.Item(6, 3).Value = "Avaliação desempenho"
$ws.Cells.Item(7, 3).Value = "Avaliação executivos"
$ws.Cells.Item(8, 3).Value = "Bônus"
$ws.Cells.Item(9, 3).Value = "Cibersegurança"
$ws.Cells.Item(10, 3).Value = "Conformidade"
$ws.Cells.Item(11, 3).Value = "Compliance"
$ws.Cells.Item(12, 3).Value = "Condenação"
$ws.Cells.Item(13, 3).Value = "Controles Internos"
$ws.Cells.Item(14, 3).Value = "Corrupção"
$ws.Cells.Item(15, 3).Value = "Crimes financeiros"
$ws.Cells.Item(16, 3).Value = "Cultura corporativa"
$ws.Cells.Item(17, 3).Value = "Desvio dinheiro"
$ws.Cells.Item(18, 3).Value = "Erros contábeis"
$ws.Cells.Item(19, 3).Value = "Escândalo"
$ws.Cells.Item(20, 3).Value = "Ética"
$ws.Cells.Item(21, 3).Value = "Fraude"
$ws.Cells.Item(22, 3).Value = "Gerenciamento Crise"
$ws.Cells.Item(23, 3).Value = "Gestão riscos"
$ws.Cells.Item(24, 3).Value = "Governança"
$ws.Cells.Item(25, 3).Value = "Honestidade"
$ws.Cells.Item(26, 3).Value = "Incidente cibernético"
$ws.Cells.Item(27, 3).Value = "Investigação"
$ws.Cells.Item(28, 3).Value = "Irregularidades"
$ws.Cells.Item(29, 3).Value = "Lavagem dinheiro"
$ws.Cells.Item(30, 3).Value = "LGPD"
$ws.Cells.Item(31, 3).Value = "Partido político"
$ws.Cells.Item(32, 3).Value = "Prejuízo investidores"
$ws.Cells.Item(33, 3).Value = "Preso"
$ws.Cells.Item(34, 3).Value = "Privacidade"
$ws.Cells.Item(35, 3).Value = "Propina"
$ws.Cells.Item(36, 3).Value = "Rombo contábil"
$ws.Cells.Item(37, 3).Value = "Segurança cibernética"
$ws.Cells.Item(38, 3).Value = "Segurança Dados"
$ws.Cells.Item(39, 3).Value = "Stakeholders"
$ws.Cells.Item(40, 3).Value = "Sonegação"
$ws.Cells.Item(41, 3).Value = "Transparência"
$ws.Cells.Item(42, 3).Value = "Vazamento dados"

# Update selection to match final view state
$ws.Range("B22").Select()
